# Meeting Recording & Burndown Update
# Sprint 3 burndown chart: record actual hours spent today on
# "Upper Floor Room F" (row 10) and "Door Locked UI" (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 3 (column E) actual hours burned on these two tasks
$ws.Range("E10").Value = 2
$ws.Range("E15").Value = 0.2

# Make sure the Actual Burndown row (and the chart feeding off it)
# reflects today's entries
$excel.CalculateFullRebuild() | Out-Null

# Clear the lingering cell selection left over from editing, back to A1
$ws.Range("A1").Select() | Out-Null
